# Update admin credentials to mohitg2/graingerlibrary across all files
# and regenerate presentation.
#
# Adds a footer textbox with the admin/project credentials to the title
# slide, and a "use only if live demo fails" note to the live-demo slide.

function RGBVal($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# --- Slide 1 (title slide): admin credentials footer ------------------
$s1 = $p.Slides.Item(1)

$shp1 = $s1.Shapes.AddTextbox(1, 72, 540, 1008, 72)
$shp1.Name = "TextBox 4"
$shp1.Fill.Visible = $false

$tf1 = $shp1.TextFrame
$tf1.WordWrap = $false
$tf1.AutoSize = 1

$tr1 = $tf1.TextRange
$tr1.Text = "Alikhan`rINFO Semester Project " + [char]0x2022 + " December 2025`r`rAdmin: mohitg2 / graingerlibrary`rgithub.com/AlikhanIllini/Final_Project_Alikhan_alikhan4"

# Only the first line ("Alikhan") carries explicit formatting in the
# source deck; the remaining lines inherit the default text style.
$firstLine1 = $tr1.Paragraphs(1, 1)
$firstLine1.Font.Size = 14
$firstLine1.Font.Color.RGB = RGBVal 0x1F 0x29 0x37
$firstLine1.ParagraphFormat.Alignment = 2

# Pin down the exact position/size (AutoSize recalculates height as text
# is typed, so re-apply the target geometry last).
$shp1.Left = 72
$shp1.Top = 540
$shp1.Width = 1008
$shp1.Height = 72

# --- Slide 6 (live demonstration): "use only if demo fails" note ------
$s6 = $p.Slides.Item(6)

$shp6 = $s6.Shapes.AddTextbox(1, 144, 540, 864, 57.600001)
$shp6.Name = "TextBox 3"
$shp6.Fill.Visible = $false

$tf6 = $shp6.TextFrame
$tf6.WordWrap = $false
$tf6.AutoSize = 1

$tr6 = $tf6.TextRange
$tr6.Text = "(Use this slide only if live demo fails - otherwise skip to live demo)"
$tr6.Font.Size = 18
$tr6.Font.Italic = $true
$tr6.Font.Color.RGB = RGBVal 0xEF 0x44 0x44
$tr6.ParagraphFormat.Alignment = 2

$shp6.Left = 144
$shp6.Top = 540
$shp6.Width = 864
$shp6.Height = 57.600001
